# -----------------------------------------------------------------------
# Applies two content changes to the "CORE COMPETENCIES" / skills section
# of the resume, per the commit's target diff:
#
# 1) Condense the three detailed "CORE COMPETENCIES" bullet paragraphs
#    into a single summary line (bare category names joined by bullets).
#
# 2) Append a new "TECHNICAL SKILLS" section (Heading 2 + three detail
#    paragraphs) right before the closing "For a more detailed..." line.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument
$bullet = [char]0x2022

# --- Change 1: collapse CORE COMPETENCIES detail paragraphs ------------
# Locate the three long paragraphs that immediately follow the
# "CORE COMPETENCIES" heading and replace them with one short line.
$coreHeadingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq "CORE COMPETENCIES") {
        $coreHeadingIndex = $i
        break
    }
}

if ($coreHeadingIndex -eq -1) {
    throw "Could not find the 'CORE COMPETENCIES' heading paragraph."
}

$firstDetail = $coreHeadingIndex + 1

# Overwrite the first detail paragraph's text with the condensed summary.
$summary = "Statistical Analysis & Machine Learning " + $bullet + `
    " Big Data & Data Engineering " + $bullet + " Data Visualization & Reporting"
$d.Paragraphs.Item($firstDetail).Range.Text = $summary

# The next two paragraphs (the old "Big Data & Data Engineering: ..." and
# "Data Visualization & Reporting: ..." detail lines) are now redundant;
# delete them outright (text + paragraph mark).
$d.Paragraphs.Item($firstDetail + 1).Range.Delete()
$d.Paragraphs.Item($firstDetail + 1).Range.Delete()

# --- Change 2: append a new TECHNICAL SKILLS section --------------------
# Find the closing "For a more detailed..." paragraph; the new section is
# inserted immediately before it.
$closingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("For a more detailed")) {
        $closingIndex = $i
        break
    }
}

if ($closingIndex -eq -1) {
    throw "Could not find the closing 'For a more detailed...' paragraph."
}

$closingRange = $d.Paragraphs.Item($closingIndex).Range
$closingRange.Collapse(1)   # wdCollapseStart

$closingRange.InsertParagraphBefore()
$closingRange.InsertParagraphBefore()
$closingRange.InsertParagraphBefore()
$closingRange.InsertParagraphBefore()

$skillsHeading = $d.Paragraphs.Item($closingIndex)
$skillsHeading.Range.Text = "TECHNICAL SKILLS"
$skillsHeading.Style = "Heading 2"

$d.Paragraphs.Item($closingIndex + 1).Range.Text = `
    "STATISTICAL ANALYSIS & MACHINE LEARNING Advanced Statistical Modeling; Predictive Analytics; Data Mining; Machine Learning"
$d.Paragraphs.Item($closingIndex + 2).Range.Text = `
    "BIG DATA & DATA ENGINEERING Big Data Processing; Data Warehousing; Cloud Platforms; Data Pipeline Optimization"
$d.Paragraphs.Item($closingIndex + 3).Range.Text = `
    "DATA VISUALIZATION & REPORTING Data Visualization; Geospatial Analysis; Interactive Dashboards; Business Intelligence"

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
